# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (A) used emoji glyphs as status codes. Swap them for
# simpler / more broadly-renderable markers:
#   📕 -> "-3"
#   📘 -> "⚠️"
#   📗 -> "✅"
#   📙 -> "+3"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ⚠️ / ✅ are plain text already, so a simple Find&Replace across the sheet
# is enough and keeps everything else (styles, types) untouched.
$ws.Cells.Replace("📘", "⚠️", $true)
$ws.Cells.Replace("📗", "✅", $true)

# "-3" / "+3" look like numbers, so Excel would silently turn them into
# numeric values if we just wrote them into a General-formatted cell.
# Force the cell to Text first so the replacement keeps its original
# (text) type, then drop the formatting change once the text is in place.
$lastRow = $ws.Cells(($ws.Rows.Count), 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()

    if ($v -eq "📕") {
        $cell.NumberFormat = "@"
        $cell.Value = "-3"
        $cell.Style = "Normal"
    } elseif ($v -eq "📙") {
        $cell.NumberFormat = "@"
        $cell.Value = "+3"
        $cell.Style = "Normal"
    }
}
